# Consolidate Cargo Type & Vehicle Type on the EoDfVUwFC sheet:
# expand the single "Vehicle" list (LDVs/HDVs/aircraft/rail/ships/motorbikes)
# into separate "passenger" and "freight" rows, and rename the header.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EoDfVUwFC")

# Widen column A slightly to fit the longer labels.
# (Excel's ColumnWidth COM property and the OOXML <col width> attribute are
# offset by 5/6 of a character, so back that out to land on exactly 22.)
$ws.Columns.Item(1).ColumnWidth = 22 - 5/6

# ---- Passenger rows (2-7), re-labelled, formulas unchanged -------------
$ws.Range("A2").Value = "passenger LDVs"
$ws.Range("A3").Value = "passenger HDVs"
$ws.Range("A4").Value = "passenger aircraft"
$ws.Range("A5").Value = "passenger rail"
$ws.Range("A6").Value = "passenger ships"
$ws.Range("A7").Value = "passenger motorbikes"

$ws.Range("B2").Formula = "='LDVs and HDVs'!B6"
$ws.Range("B3").Formula = "='LDVs and HDVs'!B7"
$ws.Range("B4").Formula = "=Aircraft!B7"
$ws.Range("B5").Formula = "='Rail and Ships'!B22"
$ws.Range("B6").Formula = "='Rail and Ships'!B23"
$ws.Range("B7").Formula = "=B2"

# ---- New freight rows (8-13), mirroring the passenger values -----------
$ws.Range("A8").Value = "freight LDVs"
$ws.Range("A9").Value = "freight HDVs"
$ws.Range("A10").Value = "freight aircraft"
$ws.Range("A11").Value = "freight rail"
$ws.Range("A12").Value = "freight ships"
$ws.Range("A13").Value = "freight motorbikes"

$ws.Range("B8").Formula = "=B2"
$ws.Range("B9").Formula = "=B3"
$ws.Range("B10").Formula = "=B4"
$ws.Range("B11").Formula = "=B5"
$ws.Range("B12").Formula = "=B6"
$ws.Range("B13").Formula = "=B7"

# Freight rows share the same 0.00 number format as passenger rows 2 and 4.
$ws.Range("B8:B13").NumberFormat = "0.00"

# ---- Header row (renamed last) ------------------------------------------
# B1 keeps its existing bold/right-aligned look, just losing the wrap (so
# touch it first - it mutates the existing style record rather than
# allocating a new one). A1 gets a brand-new italic/left/wrap style.
$ws.Range("B1").WrapText = $false

$ws.Range("A1").Font.Bold = $false
$ws.Range("A1").Font.Italic = $true
$ws.Range("A1").HorizontalAlignment = -4131   # xlLeft
$ws.Range("A1").WrapText = $true

$ws.Range("A1").Value = "Unit: dimensionless (elasticity)"
$ws.Range("B1").Value = "Value"
